# SalesTrend_CNH refresh: updated sales-trend figures for rsm (rows 2-40, cols B:J)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 84.15000000000001
$ws.Cells.Item(2, 3).Value = 87.48
$ws.Cells.Item(2, 4).Value = 85.61
$ws.Cells.Item(2, 5).Value = 58.09
$ws.Cells.Item(2, 6).Value = 97.76000000000001
$ws.Cells.Item(2, 7).Value = 93.76000000000001
$ws.Cells.Item(2, 8).Value = 89.12
$ws.Cells.Item(2, 9).Value = 50.72
$ws.Cells.Item(2, 10).Value = 94.64

$ws.Cells.Item(3, 2).Value = 94.98
$ws.Cells.Item(3, 3).Value = 121.5
$ws.Cells.Item(3, 4).Value = 91.88
$ws.Cells.Item(3, 5).Value = 32.91
$ws.Cells.Item(3, 7).Value = 114.1
$ws.Cells.Item(3, 8).Value = 102.88
$ws.Cells.Item(3, 9).Value = 48.53
$ws.Cells.Item(3, 10).Value = 121.11

$ws.Cells.Item(4, 2).Value = 90.11
$ws.Cells.Item(4, 3).Value = 93.47
$ws.Cells.Item(4, 4).Value = 91.76000000000001
$ws.Cells.Item(4, 5).Value = 9.76
$ws.Cells.Item(4, 6).Value = 53.33
$ws.Cells.Item(4, 7).Value = 79.43000000000001
$ws.Cells.Item(4, 8).Value = 114.23
$ws.Cells.Item(4, 9).Value = 49.77
$ws.Cells.Item(4, 10).Value = 78.43000000000001

$ws.Cells.Item(5, 2).Value = 68.64
$ws.Cells.Item(5, 3).Value = 123.94
$ws.Cells.Item(5, 4).Value = 76.26000000000001
$ws.Cells.Item(5, 5).Value = 41.46
$ws.Cells.Item(5, 6).Value = 97.78
$ws.Cells.Item(5, 7).Value = 172.97
$ws.Cells.Item(5, 8).Value = 90
$ws.Cells.Item(5, 9).Value = 36.56
$ws.Cells.Item(5, 10).Value = 85.88

$ws.Cells.Item(6, 2).Value = 130.2
$ws.Cells.Item(6, 3).Value = 146.11
$ws.Cells.Item(6, 4).Value = 116.83
$ws.Cells.Item(6, 5).Value = 50.79
$ws.Cells.Item(6, 6).Value = 133.33
$ws.Cells.Item(6, 7).Value = 98.04000000000001
$ws.Cells.Item(6, 8).Value = 162.69
$ws.Cells.Item(6, 9).Value = 43.12
$ws.Cells.Item(6, 10).Value = 179.31

$ws.Cells.Item(7, 2).Value = 131.49
$ws.Cells.Item(7, 3).Value = 113.53
$ws.Cells.Item(7, 4).Value = 63.71
$ws.Cells.Item(7, 5).Value = 46.26
$ws.Cells.Item(7, 6).Value = 61.54
$ws.Cells.Item(7, 7).Value = 105.43
$ws.Cells.Item(7, 8).Value = 78.83
$ws.Cells.Item(7, 9).Value = 74.66
$ws.Cells.Item(7, 10).Value = 146.03

$ws.Cells.Item(8, 2).Value = 94.65000000000001
$ws.Cells.Item(8, 3).Value = 133.6
$ws.Cells.Item(8, 4).Value = 108.11
$ws.Cells.Item(8, 5).Value = 19.51
$ws.Cells.Item(8, 6).Value = 118.52
$ws.Cells.Item(8, 7).Value = 121.99
$ws.Cells.Item(8, 8).Value = 88.22
$ws.Cells.Item(8, 9).Value = 49.13
$ws.Cells.Item(8, 10).Value = 152.78

$ws.Cells.Item(9, 2).Value = 67.65000000000001
$ws.Cells.Item(9, 3).Value = 72.31
$ws.Cells.Item(9, 4).Value = 79.2
$ws.Cells.Item(9, 5).Value = 10.4
$ws.Cells.Item(9, 6).Value = 74.45999999999999
$ws.Cells.Item(9, 7).Value = 70.75
$ws.Cells.Item(9, 8).Value = 62.58
$ws.Cells.Item(9, 9).Value = 45.58
$ws.Cells.Item(9, 10).Value = 90.25

$ws.Cells.Item(10, 2).Value = 118.24
$ws.Cells.Item(10, 3).Value = 116.86
$ws.Cells.Item(10, 4).Value = 102.97
$ws.Cells.Item(10, 5).Value = 31.37
$ws.Cells.Item(10, 6).Value = 91.23
$ws.Cells.Item(10, 7).Value = 71.5
$ws.Cells.Item(10, 8).Value = 89.19
$ws.Cells.Item(10, 9).Value = 50.19
$ws.Cells.Item(10, 10).Value = 152.73

$ws.Cells.Item(11, 2).Value = 45.63
$ws.Cells.Item(11, 3).Value = 64.65000000000001
$ws.Cells.Item(11, 4).Value = 63
$ws.Cells.Item(11, 6).Value = 26.67
$ws.Cells.Item(11, 7).Value = 75.68000000000001
$ws.Cells.Item(11, 8).Value = 48.69
$ws.Cells.Item(11, 9).Value = 28.99
$ws.Cells.Item(11, 10).Value = 64.76000000000001

$ws.Cells.Item(12, 2).Value = 90.92
$ws.Cells.Item(12, 3).Value = 73
$ws.Cells.Item(12, 4).Value = 109.06
$ws.Cells.Item(12, 5).Value = 9.98
$ws.Cells.Item(12, 6).Value = 94.12
$ws.Cells.Item(12, 7).Value = 93.98999999999999
$ws.Cells.Item(12, 8).Value = 58.11
$ws.Cells.Item(12, 9).Value = 65.61
$ws.Cells.Item(12, 10).Value = 78.79000000000001

$ws.Cells.Item(13, 2).Value = 39.76
$ws.Cells.Item(13, 3).Value = 38.53
$ws.Cells.Item(13, 4).Value = 65.56999999999999
$ws.Cells.Item(13, 5).Value = 4.88
$ws.Cells.Item(13, 6).Value = 106.67
$ws.Cells.Item(13, 7).Value = 69.56999999999999
$ws.Cells.Item(13, 8).Value = 48.3
$ws.Cells.Item(13, 9).Value = 66.27
$ws.Cells.Item(13, 10).Value = 160

$ws.Cells.Item(14, 2).Value = 30.26
$ws.Cells.Item(14, 3).Value = 35.83
$ws.Cells.Item(14, 4).Value = 30.59
$ws.Cells.Item(14, 5).Value = 13.26
$ws.Cells.Item(14, 6).Value = 58.33
$ws.Cells.Item(14, 7).Value = 41.21
$ws.Cells.Item(14, 8).Value = 50.13
$ws.Cells.Item(14, 9).Value = 19.78
$ws.Cells.Item(14, 10).Value = 30.19

$ws.Cells.Item(15, 2).Value = 121.43
$ws.Cells.Item(15, 3).Value = 135.9
$ws.Cells.Item(15, 4).Value = 111
$ws.Cells.Item(15, 5).Value = 89.56999999999999
$ws.Cells.Item(15, 6).Value = 110.39
$ws.Cells.Item(15, 7).Value = 119.17
$ws.Cells.Item(15, 8).Value = 95.93000000000001
$ws.Cells.Item(15, 9).Value = 53
$ws.Cells.Item(15, 10).Value = 105.14

$ws.Cells.Item(16, 2).Value = 108.51
$ws.Cells.Item(16, 3).Value = 141.96
$ws.Cells.Item(16, 4).Value = 111.79
$ws.Cells.Item(16, 5).Value = 65.84999999999999
$ws.Cells.Item(16, 6).Value = 141.18
$ws.Cells.Item(16, 7).Value = 71.95999999999999
$ws.Cells.Item(16, 8).Value = 74.48999999999999
$ws.Cells.Item(16, 9).Value = 32.67
$ws.Cells.Item(16, 10).Value = 105.1

$ws.Cells.Item(17, 2).Value = 127.16
$ws.Cells.Item(17, 3).Value = 121.73
$ws.Cells.Item(17, 4).Value = 137.18
$ws.Cells.Item(17, 5).Value = 121.14
$ws.Cells.Item(17, 6).Value = 100
$ws.Cells.Item(17, 7).Value = 115.66
$ws.Cells.Item(17, 8).Value = 109.64
$ws.Cells.Item(17, 9).Value = 65.01000000000001
$ws.Cells.Item(17, 10).Value = 101.9

$ws.Cells.Item(18, 2).Value = 108.94
$ws.Cells.Item(18, 3).Value = 178.5
$ws.Cells.Item(18, 4).Value = 101.29
$ws.Cells.Item(18, 5).Value = 102.49
$ws.Cells.Item(18, 6).Value = 104.76
$ws.Cells.Item(18, 7).Value = 150.48
$ws.Cells.Item(18, 8).Value = 101.38
$ws.Cells.Item(18, 9).Value = 69.36
$ws.Cells.Item(18, 10).Value = 96.59

$ws.Cells.Item(19, 2).Value = 143.38
$ws.Cells.Item(19, 3).Value = 124.09
$ws.Cells.Item(19, 4).Value = 133.88
$ws.Cells.Item(19, 5).Value = 97.05
$ws.Cells.Item(19, 6).Value = 118.52
$ws.Cells.Item(19, 7).Value = 172.55
$ws.Cells.Item(19, 8).Value = 106.17
$ws.Cells.Item(19, 9).Value = 53.13
$ws.Cells.Item(19, 10).Value = 128.79

$ws.Cells.Item(20, 2).Value = 126.23
$ws.Cells.Item(20, 3).Value = 104.67
$ws.Cells.Item(20, 4).Value = 82.67
$ws.Cells.Item(20, 5).Value = 63.41
$ws.Cells.Item(20, 6).Value = 97.78
$ws.Cells.Item(20, 7).Value = 72.73
$ws.Cells.Item(20, 8).Value = 81.52
$ws.Cells.Item(20, 9).Value = 37.03
$ws.Cells.Item(20, 10).Value = 97.98999999999999

$ws.Cells.Item(21, 2).Value = 91.45999999999999
$ws.Cells.Item(21, 3).Value = 84.52
$ws.Cells.Item(21, 4).Value = 90.84
$ws.Cells.Item(21, 5).Value = 135.18
$ws.Cells.Item(21, 6).Value = 77.88
$ws.Cells.Item(21, 7).Value = 74.06999999999999
$ws.Cells.Item(21, 8).Value = 85.23999999999999
$ws.Cells.Item(21, 9).Value = 40.75
$ws.Cells.Item(21, 10).Value = 85.83

$ws.Cells.Item(22, 2).Value = 124.85
$ws.Cells.Item(22, 3).Value = 119.11
$ws.Cells.Item(22, 4).Value = 97.03
$ws.Cells.Item(22, 5).Value = 342.91
$ws.Cells.Item(22, 6).Value = 102.22
$ws.Cells.Item(22, 7).Value = 81.48
$ws.Cells.Item(22, 8).Value = 107.29
$ws.Cells.Item(22, 9).Value = 47.14
$ws.Cells.Item(22, 10).Value = 91.59

$ws.Cells.Item(23, 2).Value = 104.56
$ws.Cells.Item(23, 3).Value = 68.76000000000001
$ws.Cells.Item(23, 4).Value = 92.59999999999999
$ws.Cells.Item(23, 5).Value = 97.05
$ws.Cells.Item(23, 6).Value = 57.14
$ws.Cells.Item(23, 7).Value = 43.14
$ws.Cells.Item(23, 8).Value = 94.87
$ws.Cells.Item(23, 9).Value = 30.19
$ws.Cells.Item(23, 10).Value = 100

$ws.Cells.Item(24, 2).Value = 52.32
$ws.Cells.Item(24, 3).Value = 60.61
$ws.Cells.Item(24, 4).Value = 56.77
$ws.Cells.Item(24, 5).Value = 14.37
$ws.Cells.Item(24, 6).Value = 40
$ws.Cells.Item(24, 7).Value = 36.04
$ws.Cells.Item(24, 8).Value = 64.53
$ws.Cells.Item(24, 9).Value = 26.54
$ws.Cells.Item(24, 10).Value = 72.81999999999999

$ws.Cells.Item(25, 2).Value = 94.20999999999999
$ws.Cells.Item(25, 3).Value = 89.27
$ws.Cells.Item(25, 4).Value = 75.14
$ws.Cells.Item(25, 5).Value = 42.71
$ws.Cells.Item(25, 6).Value = 53.33
$ws.Cells.Item(25, 7).Value = 126.22
$ws.Cells.Item(25, 8).Value = 90.95
$ws.Cells.Item(25, 9).Value = 69.93000000000001
$ws.Cells.Item(25, 10).Value = 105.98

$ws.Cells.Item(26, 2).Value = 84.65000000000001
$ws.Cells.Item(26, 3).Value = 91.48999999999999
$ws.Cells.Item(26, 4).Value = 80.20999999999999
$ws.Cells.Item(26, 5).Value = 262.13
$ws.Cells.Item(26, 6).Value = 168.89
$ws.Cells.Item(26, 7).Value = 57.14
$ws.Cells.Item(26, 8).Value = 75.59
$ws.Cells.Item(26, 9).Value = 39.29
$ws.Cells.Item(26, 10).Value = 94.87

$ws.Cells.Item(27, 2).Value = 78.06
$ws.Cells.Item(27, 3).Value = 53.61
$ws.Cells.Item(27, 4).Value = 134.17
$ws.Cells.Item(27, 5).Value = 58.96
$ws.Cells.Item(27, 6).Value = 28.57
$ws.Cells.Item(27, 7).Value = 41.9
$ws.Cells.Item(27, 8).Value = 75.64
$ws.Cells.Item(27, 9).Value = 43.83
$ws.Cells.Item(27, 10).Value = 49.44

$ws.Cells.Item(28, 2).Value = 59.91
$ws.Cells.Item(28, 3).Value = 67.84
$ws.Cells.Item(28, 4).Value = 65.42
$ws.Cells.Item(28, 5).Value = 23.74
$ws.Cells.Item(28, 6).Value = 79.31999999999999
$ws.Cells.Item(28, 7).Value = 90.31999999999999
$ws.Cells.Item(28, 8).Value = 94.01000000000001
$ws.Cells.Item(28, 9).Value = 52.69
$ws.Cells.Item(28, 10).Value = 86.19

$ws.Cells.Item(29, 2).Value = 72.98
$ws.Cells.Item(29, 3).Value = 79.05
$ws.Cells.Item(29, 4).Value = 138.25
$ws.Cells.Item(29, 5).Value = 79.81999999999999
$ws.Cells.Item(29, 6).Value = 41.03
$ws.Cells.Item(29, 7).Value = 66.67
$ws.Cells.Item(29, 8).Value = 98.53
$ws.Cells.Item(29, 9).Value = 60.71
$ws.Cells.Item(29, 10).Value = 137.5

$ws.Cells.Item(30, 2).Value = 54.36
$ws.Cells.Item(30, 3).Value = 68.01000000000001
$ws.Cells.Item(30, 4).Value = 82.77
$ws.Cells.Item(30, 5).Value = 19.95
$ws.Cells.Item(30, 6).Value = 44.44
$ws.Cells.Item(30, 7).Value = 42.42
$ws.Cells.Item(30, 8).Value = 87.14
$ws.Cells.Item(30, 9).Value = 23.85
$ws.Cells.Item(30, 10).Value = 110.48

$ws.Cells.Item(31, 2).Value = 76.55
$ws.Cells.Item(31, 3).Value = 95.09
$ws.Cells.Item(31, 4).Value = 133.06
$ws.Cells.Item(31, 5).Value = 27.11
$ws.Cells.Item(31, 6).Value = 76.19
$ws.Cells.Item(31, 7).Value = 141.18
$ws.Cells.Item(31, 8).Value = 139.13
$ws.Cells.Item(31, 9).Value = 84.66
$ws.Cells.Item(31, 10).Value = 205.71

$ws.Cells.Item(32, 2).Value = 31.15
$ws.Cells.Item(32, 3).Value = 43.42
$ws.Cells.Item(32, 4).Value = 30.24
$ws.Cells.Item(32, 5).Value = 9.859999999999999
$ws.Cells.Item(32, 6).Value = 41.03
$ws.Cells.Item(32, 7).Value = 83.33
$ws.Cells.Item(32, 8).Value = 58.74
$ws.Cells.Item(32, 9).Value = 46.64
$ws.Cells.Item(32, 10).Value = 21.88

$ws.Cells.Item(33, 2).Value = 83.25
$ws.Cells.Item(33, 3).Value = 100.62
$ws.Cells.Item(33, 4).Value = 96.44
$ws.Cells.Item(33, 5).Value = 30.84
$ws.Cells.Item(33, 6).Value = 296.3
$ws.Cells.Item(33, 7).Value = 117.65
$ws.Cells.Item(33, 8).Value = 102.44
$ws.Cells.Item(33, 9).Value = 63.36
$ws.Cells.Item(33, 10).Value = 126.32

$ws.Cells.Item(34, 2).Value = 108.94
$ws.Cells.Item(34, 3).Value = 84.43000000000001
$ws.Cells.Item(34, 4).Value = 133.59
$ws.Cells.Item(34, 5).Value = 12.31
$ws.Cells.Item(34, 6).Value = 66.67
$ws.Cells.Item(34, 7).Value = 98.55
$ws.Cells.Item(34, 8).Value = 126.06
$ws.Cells.Item(34, 9).Value = 45.01
$ws.Cells.Item(34, 10).Value = 229.33

$ws.Cells.Item(35, 2).Value = 62.26
$ws.Cells.Item(35, 3).Value = 68.98
$ws.Cells.Item(35, 4).Value = 91.37
$ws.Cells.Item(35, 5).Value = 35.79
$ws.Cells.Item(35, 6).Value = 162.16
$ws.Cells.Item(35, 7).Value = 85.88
$ws.Cells.Item(35, 8).Value = 100.23
$ws.Cells.Item(35, 9).Value = 61.03
$ws.Cells.Item(35, 10).Value = 86.84

$ws.Cells.Item(36, 2).Value = 48.96
$ws.Cells.Item(36, 3).Value = 64.48999999999999
$ws.Cells.Item(36, 4).Value = 99.67
$ws.Cells.Item(36, 5).Value = 13.82
$ws.Cells.Item(36, 6).Value = 53.33
$ws.Cells.Item(36, 7).Value = 96.55
$ws.Cells.Item(36, 8).Value = 144.63
$ws.Cells.Item(36, 9).Value = 92.56999999999999
$ws.Cells.Item(36, 10).Value = 145.83

$ws.Cells.Item(37, 2).Value = 75.31
$ws.Cells.Item(37, 3).Value = 72.36
$ws.Cells.Item(37, 4).Value = 91.33
$ws.Cells.Item(37, 6).Value = 125.49
$ws.Cells.Item(37, 7).Value = 66.67
$ws.Cells.Item(37, 8).Value = 89.43000000000001
$ws.Cells.Item(37, 9).Value = 63.84
$ws.Cells.Item(37, 10).Value = 107.94

$ws.Cells.Item(38, 2).Value = 41.56
$ws.Cells.Item(38, 3).Value = 69.14
$ws.Cells.Item(38, 4).Value = 82.61
$ws.Cells.Item(38, 5).Value = 52.61
$ws.Cells.Item(38, 6).Value = 120.63
$ws.Cells.Item(38, 7).Value = 66.67
$ws.Cells.Item(38, 8).Value = 82.97
$ws.Cells.Item(38, 9).Value = 60.14
$ws.Cells.Item(38, 10).Value = 71.11

$ws.Cells.Item(39, 2).Value = 67.61
$ws.Cells.Item(39, 3).Value = 71.37
$ws.Cells.Item(39, 4).Value = 93.55
$ws.Cells.Item(39, 5).Value = 95.94
$ws.Cells.Item(39, 6).Value = 238.1
$ws.Cells.Item(39, 7).Value = 125.71
$ws.Cells.Item(39, 8).Value = 91.83
$ws.Cells.Item(39, 9).Value = 53.91
$ws.Cells.Item(39, 10).Value = 68.56999999999999

$ws.Cells.Item(40, 2).Value = 68.36
$ws.Cells.Item(40, 3).Value = 67.45
$ws.Cells.Item(40, 4).Value = 91.89
$ws.Cells.Item(40, 5).Value = 16.33
$ws.Cells.Item(40, 6).Value = 288.89
$ws.Cells.Item(40, 7).Value = 77.52
$ws.Cells.Item(40, 8).Value = 103.88
$ws.Cells.Item(40, 10).Value = 52.03

Write-Output "SalesTrend CNH refresh applied: 347 cells updated"
